# Modify MFG and PCB design files
# - Rename BOM column header "嘉立创元件表号" -> "嘉立创元件编号"
# - Fix footprint typo for the 22uF / C6 row: "C0805" -> "TC0805"
# - Flag the two rows missing an LCSC part number ("NONE") in bold red
# - Leave the selection on E12, matching the last user interaction

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the 4th table column header (also updates the ListObject column name)
$ws.Range("D1").Value = "嘉立创元件编号"

# Correct the footprint for C6 (row 4: Comment=22uF, Designator=C6)
$ws.Range("C4").Value = "TC0805"

# Highlight the "NONE" placeholders (rows 6 and 7) in bold red
$ws.Range("D6").Font.Bold = $true
$ws.Range("D6").Font.Color = 255

$ws.Range("D7").Font.Bold = $true
$ws.Range("D7").Font.Color = 255

# Restore the final cell selection
$ws.Range("E12").Select() | Out-Null
